$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 2526
$ws.Range("I58").Value = 435
$ws.Range("J58").Value = 6011
$ws.Range("K58").Value = 1305
$ws.Range("L58").Value = 18033
$ws.Range("M58").Value = -1155
$ws.Range("N58").Value = -18333
$ws.Range("H62").Value = 2002237.5
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 2002237.5
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H70").Value = 113858.89
$ws.Range("I70").Value = 3258.2
$ws.Range("J70").Value = 252109.75
$ws.Range("K70").Value = 9774.599999999999
$ws.Range("L70").Value = 756329.25
$ws.Range("M70").Value = -9504.599999999999
$ws.Range("N70").Value = -756869.25
$ws.Range("H73").Value = 113858.89
$ws.Range("I73").Value = 3258.2
$ws.Range("J73").Value = 252109.75
$ws.Range("K73").Value = 9774.599999999999
$ws.Range("L73").Value = 756329.25
$ws.Range("M73").Value = -8838.599999999999
$ws.Range("N73").Value = -758201.25
$ws.Range("H96").Value = 1199.909
$ws.Range("I96").Value = 728.2
$ws.Range("J96").Value = 1593
$ws.Range("K96").Value = 2184.6
$ws.Range("L96").Value = 4779
$ws.Range("M96").Value = -811.6000000000004
$ws.Range("N96").Value = -7525
$ws.Range("H97").Value = 999
$ws.Range("J97").Value = 999
$ws.Range("L97").Value = 2997
$ws.Range("N97").Value = -3989
$ws.Range("H112").Value = 4133.375
$ws.Range("J112").Value = 4509.5713
$ws.Range("L112").Value = 13528.7139
$ws.Range("N112").Value = -15744.7139
$ws.Range("H132").Value = 3382.3547
$ws.Range("I132").Value = 3368.8
$ws.Range("K132").Value = 10106.4
$ws.Range("M132").Value = -7576.400000000001
$ws.Range("H138").Value = 3336.02
$ws.Range("I138").Value = 1580.75
$ws.Range("J138").Value = 4506.2
$ws.Range("K138").Value = 4742.25
$ws.Range("L138").Value = 13518.6
$ws.Range("M138").Value = 397.75
$ws.Range("N138").Value = -23798.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H60").Value = 99999
$ws.Range("J60").Value = 99999
$ws.Range("L60").Value = 99999
$ws.Range("N60").Value = -101465
$ws.Range("H61").Value = 5889.241
$ws.Range("J61").Value = 8435.75
$ws.Range("L61").Value = 8435.75
$ws.Range("N61").Value = -8859.75
$ws.Range("H74").Value = 2759.3076
$ws.Range("I74").Value = 2207.0908
$ws.Range("J74").Value = 5796.5
$ws.Range("K74").Value = 2207.0908
$ws.Range("L74").Value = 5796.5
$ws.Range("M74").Value = -1333.0908
$ws.Range("N74").Value = -7544.5
$ws.Range("H77").Value = 2759.3076
$ws.Range("I77").Value = 2207.0908
$ws.Range("J77").Value = 5796.5
$ws.Range("K77").Value = 11035.454
$ws.Range("L77").Value = 28982.5
$ws.Range("M77").Value = -6667.454
$ws.Range("N77").Value = -37718.5
$ws.Range("H97").Value = 1720.1786
$ws.Range("I97").Value = 630.44446
$ws.Range("J97").Value = 3681.7
$ws.Range("K97").Value = 630.44446
$ws.Range("L97").Value = 3681.7
$ws.Range("M97").Value = -134.44446
$ws.Range("N97").Value = -4673.7
$ws.Range("H122").Value = 1733.2693
$ws.Range("I122").Value = 1589.381
$ws.Range("K122").Value = 4768.143
$ws.Range("M122").Value = -2318.143
$ws.Range("H136").Value = 5889.241
$ws.Range("J136").Value = 8435.75
$ws.Range("L136").Value = 25307.25
$ws.Range("N136").Value = -30407.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 4879.3335
$ws.Range("I54").Value = 4879.3335
$ws.Range("K54").Value = 4879.3335
$ws.Range("M54").Value = -4395.3335
$ws.Range("H63").Value = 90000
$ws.Range("J63").Value = 90000
$ws.Range("L63").Value = 90000
$ws.Range("N63").Value = -91372
$ws.Range("H66").Value = 90000
$ws.Range("J66").Value = 90000
$ws.Range("L66").Value = 270000
$ws.Range("N66").Value = -276864
$ws.Range("H86").Value = 5937.2
$ws.Range("I86").Value = 5847.5
$ws.Range("K86").Value = 5847.5
$ws.Range("M86").Value = -4724.5
$ws.Range("H89").Value = 5937.2
$ws.Range("I89").Value = 5847.5
$ws.Range("K89").Value = 29237.5
$ws.Range("M89").Value = -23621.5
$ws.Range("H105").Value = 2330.3
$ws.Range("I105").Value = 2367.111
$ws.Range("K105").Value = 2367.111
$ws.Range("M105").Value = -620.1109999999999
$ws.Range("H134").Value = 2942.2856
$ws.Range("I134").Value = 2937.923
$ws.Range("K134").Value = 8813.769
$ws.Range("M134").Value = -6278.769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 3500
$ws.Range("J2").Value = 3500
$ws.Range("L2").Value = 3500
$ws.Range("N2").Value = -3726
$ws.Range("H25").Value = 22
$ws.Range("I25").Value = 22
$ws.Range("K25").Value = 22
$ws.Range("M25").Value = 152
$ws.Range("H56").Value = 19999
$ws.Range("J56").Value = 19999
$ws.Range("L56").Value = 19999
$ws.Range("N56").Value = -21689
$ws.Range("H58").Value = 1505.2858
$ws.Range("I58").Value = 1307.4
$ws.Range("K58").Value = 1307.4
$ws.Range("M58").Value = -1104.4
$ws.Range("H86").Value = 22229234
$ws.Range("I86").Value = 37042852
$ws.Range("K86").Value = 37042852
$ws.Range("M86").Value = -37041729
$ws.Range("H89").Value = 22229234
$ws.Range("I89").Value = 37042852
$ws.Range("K89").Value = 185214260
$ws.Range("M89").Value = -185208644
$ws.Range("H136").Value = 1505.2858
$ws.Range("I136").Value = 1307.4
$ws.Range("K136").Value = 3922.2
$ws.Range("M136").Value = -1372.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 48155004
$ws.Range("I4").Value = 77784616
$ws.Range("K4").Value = 233353848
$ws.Range("M4").Value = -233353736
$ws.Range("H11").Value = 3575.75
$ws.Range("I11").Value = 5500.8
$ws.Range("J11").Value = 367.33334
$ws.Range("K11").Value = 16502.4
$ws.Range("L11").Value = 1102.00002
$ws.Range("M11").Value = -16362.4
$ws.Range("N11").Value = -1382.00002
$ws.Range("H121").Value = 614.2
$ws.Range("I121").Value = 293.75
$ws.Range("K121").Value = 881.25
$ws.Range("M121").Value = 428.75
$ws.Range("H140").Value = 1954.4
$ws.Range("I140").Value = 1017.3333
$ws.Range("J140").Value = 3360
$ws.Range("K140").Value = 3051.9999
$ws.Range("L140").Value = 10080
$ws.Range("M140").Value = 2128.0001
$ws.Range("N140").Value = -20440

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1543.1
$ws.Range("I102").Value = 1491.375
$ws.Range("J102").Value = 1750
$ws.Range("K102").Value = 1491.375
$ws.Range("L102").Value = 1750
$ws.Range("M102").Value = 130.625
$ws.Range("N102").Value = -4994
$ws.Range("H107").Value = 398.04
$ws.Range("I107").Value = 360.25
$ws.Range("K107").Value = 360.25
$ws.Range("M107").Value = 1559.75
$ws.Range("H132").Value = 4468.077
$ws.Range("I132").Value = 4102.3228
$ws.Range("J132").Value = 5885.375
$ws.Range("K132").Value = 12306.9684
$ws.Range("L132").Value = 17656.125
$ws.Range("M132").Value = -9776.9684
$ws.Range("N132").Value = -22716.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9076.272000000001
$ws.Range("I7").Value = 5948.95
$ws.Range("K7").Value = 5948.95
$ws.Range("M7").Value = -5836.95
$ws.Range("H46").Value = 5479.4736
$ws.Range("I46").Value = 7246.6665
$ws.Range("J46").Value = 3889
$ws.Range("K46").Value = 7246.6665
$ws.Range("L46").Value = 3889
$ws.Range("M46").Value = -7058.6665
$ws.Range("N46").Value = -4265
$ws.Range("H68").Value = 10004
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 10004
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 10004
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -11502
$ws.Range("H71").Value = 10004
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 10004
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 50020
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -57508
$ws.Range("H126").Value = 9076.272000000001
$ws.Range("I126").Value = 5948.95
$ws.Range("K126").Value = 17846.85
$ws.Range("M126").Value = -15376.85
$ws.Range("H132").Value = 3750.5806
$ws.Range("I132").Value = 3727.6538
$ws.Range("J132").Value = 3869.8
$ws.Range("K132").Value = 11182.9614
$ws.Range("L132").Value = 11609.4
$ws.Range("M132").Value = -8652.9614
$ws.Range("N132").Value = -16669.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5991
$ws.Range("I62").Value = 6023.3335
$ws.Range("K62").Value = 6023.3335
$ws.Range("M62").Value = -5399.3335
$ws.Range("H65").Value = 5991
$ws.Range("I65").Value = 6023.3335
$ws.Range("K65").Value = 30116.6675
$ws.Range("M65").Value = -26996.6675
$ws.Range("H122").Value = 4877.8335
$ws.Range("I122").Value = 4237.5625
$ws.Range("K122").Value = 12712.6875
$ws.Range("M122").Value = -10262.6875
$ws.Range("H132").Value = 3234.2666
$ws.Range("I132").Value = 3171.7144
$ws.Range("K132").Value = 9515.143199999999
$ws.Range("M132").Value = -6985.143199999999
